$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 59 data, matching pattern of rows 52-58
$ws.Range("C59").Value = 10
$ws.Range("D59").Value = "Working with the filter context"
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = "Understanding arbitrarily shaped filters"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = "Understanding arbitrarily shaped filters"

$ws.Range("C59").NumberFormat = "00"
$ws.Range("E59").NumberFormat = "00"
$ws.Range("G59").NumberFormat = "00"

$ws.Range("B59").Formula = '=_xlfn.CONCAT(TEXT(C59,"00"),TEXT(E59,"00"),TEXT(G59,"00"))'

# Update selection to match new state
$ws.Range("H60").Select()
